$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SamplesTab query (row 3, column B) to remove the Tumor / Analyte Type columns
$newQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND gi.library_strategy = 'WXS'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newQuery

# Adjust the view: top-left visible cell and the active selection
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C3").Select()

$wb.Save()
